$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Аркуш1")

# Add the new value for the agricultural production metric in E13
$ws.Range("E13").Value = "Продук. сельхоз. - agrprod (тыс. руб) (8007010)"

# Update the active cell selection to reflect the new cursor position
$ws.Range("E25").Select()
